$d = $word.ActiveDocument
$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) First (empty) paragraph in the document gets a subscript vertical
#    alignment on its paragraph mark run properties.
# ---------------------------------------------------------------------------
$firstParaRange = $d.Paragraphs.Item(1).Range
$firstParaXml = '<w:p ' + $wordNs + '><w:pPr><w:rPr><w:vertAlign w:val="subscript"/></w:rPr></w:pPr></w:p>'
$null = $firstParaRange.InsertXML($firstParaXml)

# ---------------------------------------------------------------------------
# 2) Add a new row "20 | Allows adding classes" to the end of the first
#    table (Permissions table).
# ---------------------------------------------------------------------------
$permTable = $d.Tables.Item(1)
$newRow = $permTable.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "20"
$newRow.Cells.Item(2).Range.Text = "Allows adding classes"

# ---------------------------------------------------------------------------
# 3) After the second table (Supported Types), insert a new "Commissions"
#    paragraph followed by a new 2-column, 5-row table.
# ---------------------------------------------------------------------------
$typesTable = $d.Tables.Item(2)

# Create a brand new blank paragraph right after the table (pushes the
# existing trailing empty paragraphs further down) ...
$afterTypesTable = $d.Range($typesTable.Range.End, $typesTable.Range.End)
$null = $afterTypesTable.InsertParagraphBefore()

# ... then fill that new blank paragraph in with the "Commissions" text.
$commissionsParaRange = $d.Range($typesTable.Range.End, $typesTable.Range.End + 1)
$commissionsXml = '<w:p ' + $wordNs + '><w:r><w:t>Commissions</w:t></w:r></w:p>'
$null = $commissionsParaRange.InsertXML($commissionsXml)

# Insert the new table right after the "Commissions" paragraph (i.e. right
# before the first of the two original trailing empty paragraphs).
$tableInsertPoint = $d.Range($commissionsParaRange.End, $commissionsParaRange.End)
$tableXml = @"
<w:tbl $wordNs>
  <w:tblPr>
    <w:tblStyle w:val="TableGrid"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="4675"/>
    <w:gridCol w:w="4675"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>1</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:t>Presedent</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>2</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Location</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>3</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Image</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>4</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t xml:space="preserve">Image – </w:t>
        </w:r>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:t>Locaion</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>5</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>QR Code</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
"@
$null = $tableInsertPoint.InsertXML($tableXml)
